$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.445.63'
$ws.Range("E2").Value = '  -5.47%  '

$ws.Range("D3").Value = '1.837.55'
$ws.Range("E3").Value = '  -4.39%  '

$ws.Range("E4").Value = '  -0.41%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '313.39'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.68%  '

$ws.Range("E6").Value = '  -0.24%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4212'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -8.11%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3626'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -4.90%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '44.13'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -3.46%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.07211'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -6.93%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.8998'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -7.96%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '20.48'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -9.40%  '

$ws.Range("D13").Value = '1.803.94'
$ws.Range("E13").Value = '  -6.24%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.569'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -5.45%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.315'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -6.88%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.06797'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -3.00%  '

$ws.Range("E17").Value = '  -0.42%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '77.16'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -8.84%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000008977'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -5.39%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.39%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '15.30'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -8.24%  '

$ws.Range("D22").Value = '27.449.12'
$ws.Range("E22").Value = '  -5.52%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.923'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -7.95%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '10.55'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -4.37%  '

$ws.Range("D25").Value = '2.031.68'
$ws.Range("E25").Value = '  -5.49%  '

$ws.Range("E26").Value = '  -1.56%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '152.45'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -3.34%  '

$ws.Range("E28").Value = '  -4.90%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.214'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -6.83%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '110.40'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -6.30%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.662'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -9.16%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.08836'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -5.18%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.7721'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -9.89%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.492'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -11.61%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.892'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -4.05%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.24%  '

$ws.Range("E37").Value = '  -14.01%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.05342'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -5.78%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.082'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -5.83%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.01924'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -5.64%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.941'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -5.50%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '6.818'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -7.87%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.5041'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -8.15%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.1626'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -7.28%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.06605'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -4.57%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '8.185'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -12.48%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.4707'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -9.00%  '

$ws.Range("E48").Value = '  -4.96%  '

$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '10.15'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -9.18%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.625'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -7.43%  '
